$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Wrap the SPROuT logo picture (the first "First Paragraph"-styled
#    paragraph that holds an inline picture) in a one-cell table, and add an
#    empty "Image Caption"-styled paragraph below the picture inside that
#    same cell.
# ---------------------------------------------------------------------------

$logoParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.InlineShapes.Count -gt 0) {
        $logoParaIndex = $i
        break
    }
}

if ($logoParaIndex -gt 0) {
    $logoRange = $d.Paragraphs($logoParaIndex).Range

    $tableXml = '<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing">' +
        '<w:tblPr>' +
            '<w:tblStyle w:val="Table"/>' +
            '<w:tblW w:type="pct" w:w="5000"/>' +
            '<w:tblLook w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:noHBand="0" w:noVBand="0" w:val="0000"/>' +
            '<w:jc w:val="start"/>' +
        '</w:tblPr>' +
        '<w:tblGrid><w:gridCol w:w="7920"/></w:tblGrid>' +
        '<w:tr>' +
            '<w:tc>' +
                '<w:tcPr/>' +
                '<w:p>' +
                    '<w:pPr><w:jc w:val="center"/></w:pPr>' +
                    '<w:r>' +
                        '<w:drawing>' +
                            '<wp:inline>' +
                                '<wp:extent cx="2857500" cy="2857500"/>' +
                                '<wp:effectExtent b="0" l="0" r="0" t="0"/>' +
                                '<wp:docPr descr="" title="" id="21" name="Picture"/>' +
                                '<a:graphic>' +
                                    '<a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
                                        '<pic:pic>' +
                                            '<pic:nvPicPr>' +
                                                '<pic:cNvPr descr="sprout_logo_blue.jpg" id="22" name="Picture"/>' +
                                                '<pic:cNvPicPr><a:picLocks noChangeArrowheads="1" noChangeAspect="1"/></pic:cNvPicPr>' +
                                            '</pic:nvPicPr>' +
                                            '<pic:blipFill>' +
                                                '<a:blip r:embed="rId20"/>' +
                                                '<a:stretch><a:fillRect/></a:stretch>' +
                                            '</pic:blipFill>' +
                                            '<pic:spPr bwMode="auto">' +
                                                '<a:xfrm><a:off x="0" y="0"/><a:ext cx="2857500" cy="2857500"/></a:xfrm>' +
                                                '<a:prstGeom prst="rect"><a:avLst/></a:prstGeom>' +
                                                '<a:noFill/>' +
                                                '<a:ln w="9525"><a:noFill/><a:headEnd/><a:tailEnd/></a:ln>' +
                                            '</pic:spPr>' +
                                        '</pic:pic>' +
                                    '</a:graphicData>' +
                                '</a:graphic>' +
                            '</wp:inline>' +
                        '</w:drawing>' +
                    '</w:r>' +
                '</w:p>' +
                '<w:p>' +
                    '<w:pPr><w:jc w:val="center"/></w:pPr>' +
                    '<w:pPr>' +
                        '<w:jc w:val="start"/>' +
                        '<w:spacing w:before="200"/>' +
                        '<w:pStyle w:val="ImageCaption"/>' +
                    '</w:pPr>' +
                '</w:p>' +
            '</w:tc>' +
        '</w:tr>' +
    '</w:tbl>'

    $logoRange.InsertXML($tableXml)
}

# ---------------------------------------------------------------------------
# 2) Styles cleanup:
#    - remove the now-unused "Abstract Title" paragraph style
#    - restore "Abstract" style's before-spacing to 300 twips (15pt)
#    - remove the unused "Footnote Block Text" paragraph style
# ---------------------------------------------------------------------------

$abstractTitleStyle = $d.Styles("AbstractTitle")
if ($abstractTitleStyle -ne $null) {
    $abstractTitleStyle.Delete()
}

$abstractStyle = $d.Styles("Abstract")
if ($abstractStyle -ne $null) {
    $abstractStyle.ParagraphFormat.SpaceBefore = 15
}

$footnoteBlockTextStyle = $d.Styles("FootnoteBlockText")
if ($footnoteBlockTextStyle -ne $null) {
    $footnoteBlockTextStyle.Delete()
}

Write-Host "edit applied"
